# Auto-generated edit script: update crypto price/volume table
# per commit "Updated cryptos list on Fri Oct 18 22:56:03 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "68.232.90"
$ws.Range("E2").Value2 = "  +1.15%  "

# Row 3
$ws.Range("D3").Value2 = "2.637.34"
$ws.Range("E3").Value2 = "  +0.89%  "

# Row 4
$ws.Range("E4").Value2 = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "598.94"
$ws.Range("E5").Value2 = "  +1.12%  "

# Row 6
$ws.Range("E6").Value2 = "  +2.25%  "

# Row 7
$ws.Range("E7").Value2 = "  +0.00%  "

# Row 8
$ws.Range("E8").Value2 = "  -0.19%  "

# Row 9
$ws.Range("D9").Value2 = "2.635.88"
$ws.Range("E9").Value2 = "  +0.88%  "

# Row 10
$ws.Range("E10").Value2 = "  +7.24%  "

# Row 11
$ws.Range("E11").Value2 = "  -0.63%  "

# Row 12
$ws.Range("E12").Value2 = "  +1.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.349"
$ws.Range("E13").Value2 = "  +1.41%  "

# Row 14
$ws.Range("E14").Value2 = "  +2.31%  "

# Row 15
$ws.Range("E15").Value2 = "  +2.81%  "

# Row 16
$ws.Range("D16").Value2 = "3.120.20"
$ws.Range("E16").Value2 = "  +1.25%  "

# Row 17
$ws.Range("D17").Value2 = "68.240.30"
$ws.Range("E17").Value2 = "  +1.32%  "

# Row 18
$ws.Range("D18").Value2 = "2.632.87"
$ws.Range("E18").Value2 = "  +0.76%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "11.41"
$ws.Range("E19").Value2 = "  +3.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "365.74"
$ws.Range("E20").Value2 = "  -1.60%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "7.40"
$ws.Range("E21").Value2 = "  +0.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "4.27"
$ws.Range("E22").Value2 = "  -0.40%  "

# Row 23
$ws.Range("E23").Value2 = "  -0.01%  "

# Row 24
$ws.Range("E24").Value2 = "  +3.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "73.59"
$ws.Range("E25").Value2 = "  +0.18%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.999"
$ws.Range("E26").Value2 = "  -0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "10.04"
$ws.Range("E27").Value2 = "  +1.23%  "

# Row 28
$ws.Range("D28").Value2 = "2.768.17"
$ws.Range("E28").Value2 = "  +0.92%  "

# Row 29
$ws.Range("E29").Value2 = "  +5.36%  "

# Row 30
$ws.Range("E30").Value2 = "  -0.21%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "574.02"
$ws.Range("E31").Value2 = "  -1.17%  "

# Row 32
$ws.Range("B32").Value2 = "Fetch.AI"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.42"
$ws.Range("E32").Value2 = "  +3.86%  "

# Row 33
$ws.Range("B33").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "8.00"
$ws.Range("E33").Value2 = "  +4.15%  "

# Row 34
$ws.Range("E34").Value2 = "  +2.35%  "

# Row 35
$ws.Range("E35").Value2 = "  +2.63%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.999"
$ws.Range("E36").Value2 = "  +0.01%  "

# Row 37
$ws.Range("E37").Value2 = "  +3.06%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "160.43"
$ws.Range("E38").Value2 = "  +1.68%  "

# Row 39
$ws.Range("E39").Value2 = "  +0.93%  "

# Row 40
$ws.Range("E40").Value2 = "  +3.15%  "

# Row 41
$ws.Range("E41").Value2 = "  +0.46%  "

# Row 42
$ws.Range("E42").Value2 = "  +2.71%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "17.73"
$ws.Range("E43").Value2 = "  +3.47%  "

# Row 44
$ws.Range("E44").Value2 = "  +2.76%  "

# Row 45
$ws.Range("E45").Value2 = "  +13.80%  "

# Row 46
$ws.Range("E46").Value2 = "  +0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "40.46"
$ws.Range("E47").Value2 = "  -0.51%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "156.94"
$ws.Range("E48").Value2 = "  +2.21%  "

# Row 49
$ws.Range("E49").Value2 = "  +1.20%  "

# Row 50
$ws.Range("E50").Value2 = "  +1.49%  "

# Row 51
$ws.Range("E51").Value2 = "  +2.07%  "

